$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2023 -Month 6 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100107
$ws.Cells.Item($row, 8).Value = "Otros"
$ws.Cells.Item($row, 9).Value = 100107001
$ws.Cells.Item($row, 10).Value = "Caqui"
$ws.Cells.Item($row, 11).Value = "Mankaki"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 80
$ws.Cells.Item($row, 14).Value = 17000
$ws.Cells.Item($row, 15).Value = 18000
$ws.Cells.Item($row, 16).Value = 17625
$ws.Cells.Item($row, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item($row, 18).Value = "Región del Maule"
$ws.Cells.Item($row, 19).Value = 979
$ws.Cells.Item($row, 20).Value = 18
